$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Generic Backlog")

# "Buy Resharper" (row 11, under the "Professional" section) is gone - delete
# the whole row; everything below shifts up by one.
$ws.Rows("11").Delete()

# Promote "New ToastMasters speech" (now row 10) to the section's closing
# highlight style and mark it as TODO instead of IN PROGRESS.
$a10 = $ws.Range("A10")
$a10.Interior.Color = 5296274
$a10.Font.Bold = $true

$b10 = $ws.Range("B10")
$b10.Value = "TODO"
$b10.Interior.Color = 5296274
$b10.Font.Bold = $true
$b10.HorizontalAlignment = -4152

# Promote "Another pair of Glasses" (now row 16, Personal/Household section)
# to the plain-green TODO style used by the other TODO rows in that block.
$a16 = $ws.Range("A16")
$a16.Interior.Color = 5296274
$a16.Font.Bold = $false

$b16 = $ws.Range("B16")
$b16.Interior.Color = 5296274
$b16.Font.Bold = $false
$b16.HorizontalAlignment = -4152

# Match the author's final selection.
$ws.Activate()
$ws.Range("A10:B22").Select()
